$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4: latitude 16 -> 16.03, longitude -85 -> -85.15000000000001
foreach ($r in 2..4) {
    $ws.Cells.Item($r, 3).Value = 16.03
    $ws.Cells.Item($r, 4).Value = -85.15000000000001
}

# Rows 14-18 and 22: latitude 16.39 -> 16.33, longitude -86.26000000000001 -> -86.59
foreach ($r in @(14, 15, 16, 17, 18, 22)) {
    $ws.Cells.Item($r, 3).Value = 16.33
    $ws.Cells.Item($r, 4).Value = -86.59
}
